$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10
$ws.Range('A10').Value = 'command-r'
$ws.Range('B10').Value = 'llama3:70b'
$ws.Range('C10').Value = 100
$ws.Range('D10').Value = 400
$ws.Range('E10').Value = 2709.84
$ws.Range('F10').Value = 1335.52
$ws.Range('G10').Value = 0.28125
$ws.Range('H10').Value = 'logs\command_r_llama3_70b_100_400_test_match.txt'
$ws.Range('I10').Value = 1335.52
$ws.Range('J10').Value = 0.90625
$ws.Range('K10').Value = 'logs\command_r_llama3_70b_100_400_test_correct.txt'
$ws.Range('L10').Value = 1335.52
$ws.Range('M10').Value = 0.94375
$ws.Range('N10').Value = 'logs\command_r_llama3_70b_100_400_test_executable.txt'
$ws.Range('O10').Value = 0
$ws.Range('AC10').Value = 'text'
$ws.Range('AD10').Value = 140
$ws.Range('AE10').Value = 0.1
$ws.Range('AF10').Value = 150
$ws.Range('AG10').Value = 0.9
$ws.Range('AH10').Value = 5
$ws.Range('AI10').Value = 1
$ws.Range('AJ10').Value = 1.2
$ws.Range('AK10').Value = 1
$ws.Range('AL10').Value = 1024
$ws.Range('AM10').Value = 0
$ws.Range('AN10').Value = 1374.32
$ws.Range('AO10').Value = 0.3125
$ws.Range('AP10').Value = 'logs\command_r_llama3_70b_100_400_test_fewshot_match.txt'
$ws.Range('AQ10').Value = 1374.32
$ws.Range('AR10').Value = 0.8875
$ws.Range('AS10').Value = 'logs\command_r_llama3_70b_100_400_test_fewshot_correct.txt'
$ws.Range('AT10').Value = 1374.32
$ws.Range('AU10').Value = 0.925
$ws.Range('AV10').Value = 'logs\command_r_llama3_70b_100_400_test_fewshot_executable.txt'
$ws.Range('AW10').Value = 0

# Row 11
$ws.Range('A11').Value = 'aya:35b'
$ws.Range('B11').Value = 'llama3:70b'
$ws.Range('C11').Value = 100
$ws.Range('D11').Value = 400
$ws.Range('E11').Value = 2695.8
$ws.Range('F11').Value = 1345.42
$ws.Range('G11').Value = 0.28125
$ws.Range('H11').Value = 'logs\aya_35b_llama3_70b_100_400_test_match.txt'
$ws.Range('I11').Value = 1345.42
$ws.Range('J11').Value = 0.875
$ws.Range('K11').Value = 'logs\aya_35b_llama3_70b_100_400_test_correct.txt'
$ws.Range('L11').Value = 1345.42
$ws.Range('M11').Value = 0.93125
$ws.Range('N11').Value = 'logs\aya_35b_llama3_70b_100_400_test_executable.txt'
$ws.Range('O11').Value = 0
$ws.Range('AC11').Value = 'text'
$ws.Range('AD11').Value = 140
$ws.Range('AE11').Value = 0.1
$ws.Range('AF11').Value = 150
$ws.Range('AG11').Value = 0.9
$ws.Range('AH11').Value = 5
$ws.Range('AI11').Value = 1
$ws.Range('AJ11').Value = 1.2
$ws.Range('AK11').Value = 1
$ws.Range('AL11').Value = 1024
$ws.Range('AM11').Value = 0
$ws.Range('AN11').Value = 1350.38
$ws.Range('AO11').Value = 0.2875
$ws.Range('AP11').Value = 'logs\aya_35b_llama3_70b_100_400_test_fewshot_match.txt'
$ws.Range('AQ11').Value = 1350.38
$ws.Range('AR11').Value = 0.85
$ws.Range('AS11').Value = 'logs\aya_35b_llama3_70b_100_400_test_fewshot_correct.txt'
$ws.Range('AT11').Value = 1350.38
$ws.Range('AU11').Value = 0.9375
$ws.Range('AV11').Value = 'logs\aya_35b_llama3_70b_100_400_test_fewshot_executable.txt'
$ws.Range('AW11').Value = 0

# Row 12
$ws.Range('A12').Value = 'qwen2:7b-instruct-q5_K_M'
$ws.Range('B12').Value = 'llama3:70b'
$ws.Range('C12').Value = 100
$ws.Range('D12').Value = 400
$ws.Range('E12').Value = 2632.75
$ws.Range('F12').Value = 1309.74
$ws.Range('G12').Value = 0.275
$ws.Range('H12').Value = 'logs\qwen2_7b_instruct_q5_K_M_llama3_70b_100_400_test_match.txt'
$ws.Range('I12').Value = 1309.74
$ws.Range('J12').Value = 0.9
$ws.Range('K12').Value = 'logs\qwen2_7b_instruct_q5_K_M_llama3_70b_100_400_test_correct.txt'
$ws.Range('L12').Value = 1309.74
$ws.Range('M12').Value = 0.93125
$ws.Range('N12').Value = 'logs\qwen2_7b_instruct_q5_K_M_llama3_70b_100_400_test_executable.txt'
$ws.Range('O12').Value = 0
$ws.Range('AC12').Value = 'text'
$ws.Range('AD12').Value = 140
$ws.Range('AE12').Value = 0.1
$ws.Range('AF12').Value = 150
$ws.Range('AG12').Value = 0.9
$ws.Range('AH12').Value = 5
$ws.Range('AI12').Value = 1
$ws.Range('AJ12').Value = 1.2
$ws.Range('AK12').Value = 1
$ws.Range('AL12').Value = 1024
$ws.Range('AM12').Value = 0
$ws.Range('AN12').Value = 1323.01
$ws.Range('AO12').Value = 0.3
$ws.Range('AP12').Value = 'logs\qwen2_7b_instruct_q5_K_M_llama3_70b_100_400_test_fewshot_match.txt'
$ws.Range('AQ12').Value = 1323.01
$ws.Range('AR12').Value = 0.86875
$ws.Range('AS12').Value = 'logs\qwen2_7b_instruct_q5_K_M_llama3_70b_100_400_test_fewshot_correct.txt'
$ws.Range('AT12').Value = 1323.01
$ws.Range('AU12').Value = 0.94375
$ws.Range('AV12').Value = 'logs\qwen2_7b_instruct_q5_K_M_llama3_70b_100_400_test_fewshot_executable.txt'
$ws.Range('AW12').Value = 0

# Row 13
$ws.Range('A13').Value = 'llama3:8b-instruct-fp16'
$ws.Range('B13').Value = 'llama3:70b'
$ws.Range('C13').Value = 100
$ws.Range('D13').Value = 400
$ws.Range('E13').Value = 2657.96
$ws.Range('F13').Value = 1317.85
$ws.Range('G13').Value = 0.3
$ws.Range('H13').Value = 'logs\llama3_8b_instruct_fp16_llama3_70b_100_400_test_match.txt'
$ws.Range('I13').Value = 1317.85
$ws.Range('J13').Value = 0.925
$ws.Range('K13').Value = 'logs\llama3_8b_instruct_fp16_llama3_70b_100_400_test_correct.txt'
$ws.Range('L13').Value = 1317.85
$ws.Range('M13').Value = 0.9375
$ws.Range('N13').Value = 'logs\llama3_8b_instruct_fp16_llama3_70b_100_400_test_executable.txt'
$ws.Range('O13').Value = 0
$ws.Range('AC13').Value = 'text'
$ws.Range('AD13').Value = 140
$ws.Range('AE13').Value = 0.1
$ws.Range('AF13').Value = 150
$ws.Range('AG13').Value = 0.9
$ws.Range('AH13').Value = 5
$ws.Range('AI13').Value = 1
$ws.Range('AJ13').Value = 1.2
$ws.Range('AK13').Value = 1
$ws.Range('AL13').Value = 1024
$ws.Range('AM13').Value = 0
$ws.Range('AN13').Value = 1340.11
$ws.Range('AO13').Value = 0.2625
$ws.Range('AP13').Value = 'logs\llama3_8b_instruct_fp16_llama3_70b_100_400_test_fewshot_match.txt'
$ws.Range('AQ13').Value = 1340.11
$ws.Range('AR13').Value = 0.84375
$ws.Range('AS13').Value = 'logs\llama3_8b_instruct_fp16_llama3_70b_100_400_test_fewshot_correct.txt'
$ws.Range('AT13').Value = 1340.11
$ws.Range('AU13').Value = 0.93125
$ws.Range('AV13').Value = 'logs\llama3_8b_instruct_fp16_llama3_70b_100_400_test_fewshot_executable.txt'
$ws.Range('AW13').Value = 0

# Row 14
$ws.Range('A14').Value = 'codegemma:7b-code-fp16'
$ws.Range('B14').Value = 'llama3:70b'
$ws.Range('C14').Value = 100
$ws.Range('D14').Value = 400
$ws.Range('E14').Value = 3991.24
$ws.Range('F14').Value = 1958.43
$ws.Range('G14').Value = 0.23125
$ws.Range('H14').Value = 'logs\codegemma_7b_code_fp16_llama3_70b_100_400_test_match.txt'
$ws.Range('I14').Value = 1958.43
$ws.Range('J14').Value = 0.68125
$ws.Range('K14').Value = 'logs\codegemma_7b_code_fp16_llama3_70b_100_400_test_correct.txt'
$ws.Range('L14').Value = 1958.43
$ws.Range('M14').Value = 0.425
$ws.Range('N14').Value = 'logs\codegemma_7b_code_fp16_llama3_70b_100_400_test_executable.txt'
$ws.Range('O14').Value = 0
$ws.Range('AC14').Value = 'text'
$ws.Range('AD14').Value = 140
$ws.Range('AE14').Value = 0.1
$ws.Range('AF14').Value = 150
$ws.Range('AG14').Value = 0.9
$ws.Range('AH14').Value = 5
$ws.Range('AI14').Value = 1
$ws.Range('AJ14').Value = 1.2
$ws.Range('AK14').Value = 1
$ws.Range('AL14').Value = 1024
$ws.Range('AM14').Value = 0
$ws.Range('AN14').Value = 2032.81
$ws.Range('AO14').Value = 0.33125
$ws.Range('AP14').Value = 'logs\codegemma_7b_code_fp16_llama3_70b_100_400_test_fewshot_match.txt'
$ws.Range('AQ14').Value = 2032.81
$ws.Range('AR14').Value = 0.88125
$ws.Range('AS14').Value = 'logs\codegemma_7b_code_fp16_llama3_70b_100_400_test_fewshot_correct.txt'
$ws.Range('AT14').Value = 2032.81
$ws.Range('AU14').Value = 0.25
$ws.Range('AV14').Value = 'logs\codegemma_7b_code_fp16_llama3_70b_100_400_test_fewshot_executable.txt'
$ws.Range('AW14').Value = 0
